$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate the Arabic column headers to their English equivalents.
# (This also renames the corresponding Excel Table's column names,
# since the table header cells and the ListObject column names are
# kept in sync by Excel.)
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "gender"
$ws.Range("C1").Value = "academic_rank"
$ws.Range("D1").Value = "college"
$ws.Range("E1").Value = "department"
$ws.Range("F1").Value = "research_interests"
$ws.Range("G1").Value = "phone"
$ws.Range("H1").Value = "email"
$ws.Range("I1").Value = "notes"

# Turn off the table's AutoFilter dropdown buttons.
$tbl = $ws.ListObjects.Item(1)
$tbl.ShowAutoFilter = $false

# Resize the columns to fit the new (shorter) English header text.
$ws.Columns.Item(1).ColumnWidth = 6.857142857142857
$ws.Columns.Item(2).ColumnWidth = 8.714285714285714
$ws.Columns.Item(3).ColumnWidth = 18.428571428571427
$ws.Columns.Item(4).ColumnWidth = 8.857142857142858
$ws.Columns.Item(5).ColumnWidth = 13.714285714285714
$ws.Columns.Item(6).ColumnWidth = 22.714285714285715
$ws.Columns.Item(7).ColumnWidth = 7.714285714285714
$ws.Columns.Item(8).ColumnWidth = 6.571428571428571
$ws.Columns.Item(9).ColumnWidth = 6.857142857142857

# Move the active selection to C3.
$ws.Range("C3").Select() | Out-Null
